$wb = $excel.ActiveWorkbook

# --- Sheet: Validation Metrics ---
$wsVal = $wb.Worksheets.Item("Validation Metrics")
$wsVal.Range("B2").Value = 0.4806373119354248
$wsVal.Range("B3").Value = 0.9164574146270752

# --- Sheet: Classification Report ---
$wsCls = $wb.Worksheets.Item("Classification Report")

# Row 3 (class "1")
$wsCls.Range("B3").Value = 0.5114401076716016
$wsCls.Range("C3").Value = 0.9973753280839895
$wsCls.Range("D3").Value = 0.6761565836298933
$wsCls.Range("E3").Value = 381

# Row 4 (class "2")
$wsCls.Range("B4").Value = 0.999
$wsCls.Range("D4").Value = 0.999

# Row 7 (accuracy)
$wsCls.Range("B7").Value = 0.9164574298105456
$wsCls.Range("C7").Value = 0.9164574298105456
$wsCls.Range("D7").Value = 0.9164574298105456
$wsCls.Range("E7").Value = 0.9164574298105456

# Row 8 (macro avg)
$wsCls.Range("B8").Value = 0.9016882213345202
$wsCls.Range("C8").Value = 0.9264750656167979
$wsCls.Range("D8").Value = 0.8903819470308628
$wsCls.Range("E8").Value = 4381

# Row 9 (weighted avg)
$wsCls.Range("B9").Value = 0.9568271353626751
$wsCls.Range("C9").Value = 0.9164574298105456
$wsCls.Range("D9").Value = 0.920650264754031
$wsCls.Range("E9").Value = 4381

# --- Sheet: Confusion Matrix ---
$wsConf = $wb.Worksheets.Item("Confusion Matrix")
$wsConf.Range("C3").Value = 380
$wsConf.Range("D3").Value = 1
